{"js": "// Replace the date line and each two-digit \u00f7 one-digit division fact with\n// the updated values from the next day's worksheet output.\nconst replacements = [\n  [\"2024-10-16 Wednesday\", \"2024-10-17 Thursday\"],\n  [\"24\u00f77=3, 3\", \"59\u00f74=14, 3\"],\n  [\"82\u00f78=10, 2\", \"89\u00f74=22, 1\"],\n  [\"27\u00f78=3, 3\", \"69\u00f73=23, 0\"],\n  [\"79\u00f75=15, 4\", \"17\u00f78=2, 1\"],\n  [\"54\u00f77=7, 5\", \"90\u00f73=30, 0\"],\n  [\"23\u00f76=3, 5\", \"46\u00f78=5, 6\"],\n  [\"45\u00f73=15, 0\", \"82\u00f72=41, 0\"],\n  [\"43\u00f76=7, 1\", \"95\u00f76=15, 5\"],\n  [\"43\u00f74=10, 3\", \"27\u00f79=3, 0\"],\n  [\"25\u00f76=4, 1\", \"63\u00f74=15, 3\"],\n  [\"51\u00f73=17, 0\", \"90\u00f72=45, 0\"],\n  [\"11\u00f74=2, 3\", \"78\u00f72=39, 0\"],\n  [\"10\u00f72=5, 0\", \"98\u00f76=16, 2\"],\n  [\"68\u00f72=34, 0\", \"48\u00f73=16, 0\"],\n  [\"52\u00f76=8, 4\", \"40\u00f72=20, 0\"],\n  [\"52\u00f75=10, 2\", \"94\u00f74=23, 2\"],\n  [\"49\u00f76=8, 1\", \"45\u00f77=6, 3\"],\n  [\"88\u00f75=17, 3\", \"11\u00f74=2, 3\"],\n  [\"20\u00f74=5, 0\", \"54\u00f77=7, 5\"],\n  [\"36\u00f75=7, 1\", \"57\u00f77=8, 1\"],\n  [\"24\u00f74=6, 0\", \"57\u00f78=7, 1\"],\n  [\"75\u00f79=8, 3\", \"30\u00f75=6, 0\"],\n  [\"28\u00f72=14, 0\", \"75\u00f76=12, 3\"],\n  [\"24\u00f73=8, 0\", \"58\u00f72=29, 0\"],\n  [\"95\u00f77=13, 4\", \"46\u00f74=11, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit \u00f7 one-digit division fact with\n# the updated values from the next day's worksheet output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-10-16 Wednesday\", \"2024-10-17 Thursday\"),\n  @(\"24\u00f77=3, 3\", \"59\u00f74=14, 3\"),\n  @(\"82\u00f78=10, 2\", \"89\u00f74=22, 1\"),\n  @(\"27\u00f78=3, 3\", \"69\u00f73=23, 0\"),\n  @(\"79\u00f75=15, 4\", \"17\u00f78=2, 1\"),\n  @(\"54\u00f77=7, 5\", \"90\u00f73=30, 0\"),\n  @(\"23\u00f76=3, 5\", \"46\u00f78=5, 6\"),\n  @(\"45\u00f73=15, 0\", \"82\u00f72=41, 0\"),\n  @(\"43\u00f76=7, 1\", \"95\u00f76=15, 5\"),\n  @(\"43\u00f74=10, 3\", \"27\u00f79=3, 0\"),\n  @(\"25\u00f76=4, 1\", \"63\u00f74=15, 3\"),\n  @(\"51\u00f73=17, 0\", \"90\u00f72=45, 0\"),\n  @(\"11\u00f74=2, 3\", \"78\u00f72=39, 0\"),\n  @(\"10\u00f72=5, 0\", \"98\u00f76=16, 2\"),\n  @(\"68\u00f72=34, 0\", \"48\u00f73=16, 0\"),\n  @(\"52\u00f76=8, 4\", \"40\u00f72=20, 0\"),\n  @(\"52\u00f75=10, 2\", \"94\u00f74=23, 2\"),\n  @(\"49\u00f76=8, 1\", \"45\u00f77=6, 3\"),\n  @(\"88\u00f75=17, 3\", \"11\u00f74=2, 3\"),\n  @(\"20\u00f74=5, 0\", \"54\u00f77=7, 5\"),\n  @(\"36\u00f75=7, 1\", \"57\u00f77=8, 1\"),\n  @(\"24\u00f74=6, 0\", \"57\u00f78=7, 1\"),\n  @(\"75\u00f79=8, 3\", \"30\u00f75=6, 0\"),\n  @(\"28\u00f72=14, 0\", \"75\u00f76=12, 3\"),\n  @(\"24\u00f73=8, 0\", \"58\u00f72=29, 0\"),\n  @(\"95\u00f77=13, 4\", \"46\u00f74=11, 2\")\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $found = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n\n  if (-not $found) {\n    Write-Output \"NOT FOUND: $old\"\n  }\n}\n"}
